$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "monte carlo iterations" value in B9
$ws.Range("B9").Value = 10000000

# Move the active selection to B9 (as last edited cell)
$ws.Range("B9").Select()
